$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in A2:B5 per the cluster analysis re-ordering
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 428

$ws.Range("A3").Value = 22
$ws.Range("B3").Value = 253

$ws.Range("A4").Value = 21
$ws.Range("B4").Value = 202

$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 119
